# Updates cryptos list (prices / 1h volume %, and a few re-ranked rows)
# to match the data scraped on Mon May 22 11:24:29 UTC 2023.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.995.67'
$ws.Range("E2").Value = '  -0.52%  '
$ws.Range("D3").Value = '1.829.11'
$ws.Range("E3").Value = '  +0.14%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.66'
$ws.Range("E5").Value = '  -0.40%  '
$ws.Range("E6").Value = '  -0.15%  '
$ws.Range("E7").Value = '  -1.17%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3708'
$ws.Range("E8").Value = '  +1.69%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07364'
$ws.Range("E9").Value = '  -0.45%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8755'
$ws.Range("E10").Value = '  -0.48%  '
$ws.Range("B11").Value = 'TRON'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07900'
$ws.Range("E11").Value = '  +7.71%  '
$ws.Range("B12").Value = 'Solana'
$ws.Range("C12").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '19.96'
$ws.Range("E12").Value = '  -1.87%  '
$ws.Range("D13").Value = '1.791.89'
$ws.Range("E13").Value = '  -7.11%  '
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.595'
$ws.Range("E14").Value = '  +1.25%  '
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.359'
$ws.Range("E15").Value = '  -0.39%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '92.01'
$ws.Range("E16").Value = '  -1.29%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.011'
$ws.Range("E17").Value = '  +0.27%  '
$ws.Range("E18").Value = '  +2.16%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.009'
$ws.Range("E19").Value = '  -0.19%  '
$ws.Range("E20").Value = '  +0.48%  '
$ws.Range("D21").Value = '27.195.70'
$ws.Range("E21").Value = '  -1.77%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.162'
$ws.Range("E22").Value = '  -1.54%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.58'
$ws.Range("E23").Value = '  +0.15%  '
$ws.Range("D24").Value = '1.999.29'
$ws.Range("E24").Value = '  -4.83%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.72'
$ws.Range("E25").Value = '  +0.73%  '
$ws.Range("E26").Value = '  -2.55%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.23'
$ws.Range("E27").Value = '  -1.70%  '
$ws.Range("E28").Value = '  -1.54%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.130'
$ws.Range("E29").Value = '  -0.84%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '115.40'
$ws.Range("E30").Value = '  -0.69%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08883'
$ws.Range("E31").Value = '  -0.58%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.988'
$ws.Range("E32").Value = '  +1.63%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.444'
$ws.Range("E33").Value = '  -1.32%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7264'
$ws.Range("E34").Value = '  -1.94%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.135'
$ws.Range("E35").Value = '  -2.69%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.487'
$ws.Range("E36").Value = '  +3.35%  '
$ws.Range("E37").Value = '  -0.76%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01950'
$ws.Range("E38").Value = '  +0.40%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05246'
$ws.Range("E39").Value = '  -0.92%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.316'
$ws.Range("E40").Value = '  +1.83%  '
$ws.Range("E41").Value = '  -0.24%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5192'
$ws.Range("E42").Value = '  -1.19%  '
$ws.Range("B43").Value = 'Frax'
$ws.Range("C43").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8712'
$ws.Range("E43").Value = '  -13.78%  '
$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1626'
$ws.Range("E44").Value = '  -0.92%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.213'
$ws.Range("E45").Value = '  -1.83%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4840'
$ws.Range("E46").Value = '  -0.63%  '
$ws.Range("E47").Value = '  -0.15%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.23'
$ws.Range("E48").Value = '  -1.46%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '102.92'
$ws.Range("E49").Value = '  -1.28%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.625'
$ws.Range("E50").Value = '  -1.50%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06231'
$ws.Range("E51").Value = '  -1.10%  '

